# Update the "new values" generated during a CreateEntity / SendNotify demo run.
# Four of the five test sheets (Sheet1..Sheet4) get fresh mobile numbers in the
# MobileNumber / Enquiry_PhoneNumber / Lead_PN / Sales_PN columns (G, AF, AV, AZ)
# and Sheet1 / Sheet2 additionally get a refreshed "CurrentTime" stamp in BB2.

$wb = $excel.ActiveWorkbook

$newPhoneNumbers = @{
    "G2"  = "9840025402"
    "AF2" = "9840005389"
    "AV2" = "9840061003"
    "AZ2" = "9840068012"
}

foreach ($sheetName in @("Sheet1", "Sheet2", "Sheet3", "Sheet4")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $newPhoneNumbers.Keys) {
        $ws.Range($addr).Value = $newPhoneNumbers[$addr]
    }
}

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("BB2").Value = "CT: Mon, Jan 06, 2025 at 6:37 PM"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("BB2").Value = "CT: Mon, Jan 06, 2025 at 6:50 PM"
